$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Fix merged cell ranges in column A-D for rows 8-20 -> 8-21
#    (a new row of data - "Earthing Lug" - is added at row 21,
#     and row 20 now holds the new "Equal Angle Bar" detail row)
# ---------------------------------------------------------------
$ws.Range("A8:A20").UnMerge()
$ws.Range("B8:B20").UnMerge()
$ws.Range("C8:C20").UnMerge()
$ws.Range("D8:D20").UnMerge()

$ws.Range("A8:A21").Merge()
$ws.Range("B8:B21").Merge()
$ws.Range("C8:C21").Merge()
$ws.Range("D8:D21").Merge()

# ---------------------------------------------------------------
# 2) Update cell values (row 8 equipment header + all part rows)
# ---------------------------------------------------------------
$ws.Range("C8").Value2 = "MLK_PMT_10102_-_V-002_1"
$ws.Range("D8").Value2 = "Expansion Tank"

$ws.Range("E8").Value2 = "Tube Bundle"
$ws.Range("G8").Value2 = "DMSO"
$ws.Range("H8").Value2 = "Stainless Steel"
$ws.Range("I8").Value2 = "SA-240"
$ws.Range("J8").Value2 = "316L"
$ws.Range("K8").Value2 = "N/A"
$ws.Range("L8").Value2 = "100 °C"
$ws.Range("M8").Value2 = "1.1 Bar G"
$ws.Range("N8").Value2 = "100 °C"
$ws.Range("O8").Value2 = "1.0 Bar G"
$ws.Range("E9").Value2 = "Dish Head (Top) (6mm, 2:1 Type)"
$ws.Range("G9").Value2 = "DMSO"
$ws.Range("H9").Value2 = "Stainless Steel"
$ws.Range("I9").Value2 = "SA-240"
$ws.Range("J9").Value2 = "316L"
$ws.Range("K9").Value2 = "N/A"
$ws.Range("L9").Value2 = "100 °C"
$ws.Range("M9").Value2 = "1.1 Bar G"
$ws.Range("N9").Value2 = "100 °C"
$ws.Range("O9").Value2 = "1.0 Bar G"
$ws.Range("E10").Value2 = "Dish Head (Bottom) (6mm, 2:1 Type)"
$ws.Range("G10").Value2 = "DMSO"
$ws.Range("H10").Value2 = "Stainless Steel"
$ws.Range("I10").Value2 = "SA-240"
$ws.Range("J10").Value2 = "316L"
$ws.Range("K10").Value2 = "N/A"
$ws.Range("L10").Value2 = "100 °C"
$ws.Range("M10").Value2 = "1.1 Bar G"
$ws.Range("N10").Value2 = "100 °C"
$ws.Range("O10").Value2 = "1.0 Bar G"
$ws.Range("E11").Value2 = "Seamless Pipe DN50 x 67 SCH 40s"
$ws.Range("G11").Value2 = "DMSO"
$ws.Range("H11").Value2 = "Stainless Steel"
$ws.Range("I11").Value2 = "SA-312"
$ws.Range("J11").Value2 = "TP 316L"
$ws.Range("K11").Value2 = "N/A"
$ws.Range("L11").Value2 = "100 °C"
$ws.Range("M11").Value2 = "1.1 Bar G"
$ws.Range("N11").Value2 = "100 °C"
$ws.Range("O11").Value2 = "1.0 Bar G"
$ws.Range("E12").Value2 = "Seamless Pipe DN25 x 100 SCH 40s"
$ws.Range("G12").Value2 = "DMSO"
$ws.Range("H12").Value2 = "Stainless Steel"
$ws.Range("I12").Value2 = "SA-312"
$ws.Range("J12").Value2 = "TP 316L"
$ws.Range("K12").Value2 = "N/A"
$ws.Range("L12").Value2 = "100 °C"
$ws.Range("M12").Value2 = "1.1 Bar G"
$ws.Range("N12").Value2 = "100 °C"
$ws.Range("O12").Value2 = "1.0 Bar G"
$ws.Range("E13").Value2 = "Seamless Pipe DN150 x 112 SCH 40s"
$ws.Range("G13").Value2 = "DMSO"
$ws.Range("H13").Value2 = "Stainless Steel"
$ws.Range("I13").Value2 = "SA-312"
$ws.Range("J13").Value2 = "TP 316L"
$ws.Range("K13").Value2 = "N/A"
$ws.Range("L13").Value2 = "100 °C"
$ws.Range("M13").Value2 = "1.1 Bar G"
$ws.Range("N13").Value2 = "100 °C"
$ws.Range("O13").Value2 = "1.0 Bar G"
$ws.Range("E14").Value2 = "Flange DN50 Class 150 WNRF SCH 40s"
$ws.Range("G14").Value2 = "DMSO"
$ws.Range("H14").Value2 = "Stainless Steel"
$ws.Range("I14").Value2 = "SA-182"
$ws.Range("J14").Value2 = "F316L"
$ws.Range("K14").Value2 = "N/A"
$ws.Range("L14").Value2 = "100 °C"
$ws.Range("M14").Value2 = "1.1 Bar G"
$ws.Range("N14").Value2 = "100 °C"
$ws.Range("O14").Value2 = "1.0 Bar G"
$ws.Range("E15").Value2 = "Flange DN25 Class 150 WNRF SCH 40s"
$ws.Range("G15").Value2 = "DMSO"
$ws.Range("H15").Value2 = "Stainless Steel"
$ws.Range("I15").Value2 = "SA-182"
$ws.Range("J15").Value2 = "F316L"
$ws.Range("K15").Value2 = "N/A"
$ws.Range("L15").Value2 = "100 °C"
$ws.Range("M15").Value2 = "1.1 Bar G"
$ws.Range("N15").Value2 = "100 °C"
$ws.Range("O15").Value2 = "1.0 Bar G"
$ws.Range("E16").Value2 = "Flange DN150 Class 150 WNRF SCH 40s"
$ws.Range("G16").Value2 = "DMSO"
$ws.Range("H16").Value2 = "Stainless Steel"
$ws.Range("I16").Value2 = "SA-182"
$ws.Range("J16").Value2 = "F316L"
$ws.Range("K16").Value2 = "N/A"
$ws.Range("L16").Value2 = "100 °C"
$ws.Range("M16").Value2 = "1.1 Bar G"
$ws.Range("N16").Value2 = "100 °C"
$ws.Range("O16").Value2 = "1.0 Bar G"
$ws.Range("E17").Value2 = "Base Plate (150 x 150 x 5mm Thk)"
$ws.Range("G17").Value2 = "DMSO"
$ws.Range("H17").Value2 = "Stainless Steel"
$ws.Range("I17").Value2 = "SA-240"
$ws.Range("J17").Value2 = "Gr. 304"
$ws.Range("K17").Value2 = "N/A"
$ws.Range("L17").Value2 = "100 °C"
$ws.Range("M17").Value2 = "1.1 Bar G"
$ws.Range("N17").Value2 = "100 °C"
$ws.Range("O17").Value2 = "1.0 Bar G"
$ws.Range("E18").Value2 = "Lifting Lug Plate (190 x 80 x 6mm Thk)"
$ws.Range("G18").Value2 = "DMSO"
$ws.Range("H18").Value2 = "Stainless Steel"
$ws.Range("I18").Value2 = "SA-240"
$ws.Range("J18").Value2 = "Gr. 304"
$ws.Range("K18").Value2 = "N/A"
$ws.Range("L18").Value2 = "100 °C"
$ws.Range("M18").Value2 = "1.1 Bar G"
$ws.Range("N18").Value2 = "100 °C"
$ws.Range("O18").Value2 = "1.0 Bar G"
$ws.Range("E19").Value2 = "Doubler Plate (120 x 150 x 12.7mm Thk)"
$ws.Range("G19").Value2 = "DMSO"
$ws.Range("H19").Value2 = "Stainless Steel"
$ws.Range("I19").Value2 = "SA-240"
$ws.Range("J19").Value2 = "Gr. 304"
$ws.Range("K19").Value2 = "N/A"
$ws.Range("L19").Value2 = "100 °C"
$ws.Range("M19").Value2 = "1.1 Bar G"
$ws.Range("N19").Value2 = "100 °C"
$ws.Range("O19").Value2 = "1.0 Bar G"
$ws.Range("E20").Value2 = "Equal Angle Bar (3`" x 3`" x 1/4`" Thk)"
$ws.Range("G20").Value2 = "DMSO"
$ws.Range("H20").Value2 = "Stainless Steel"
$ws.Range("I20").Value2 = "SA-240"
$ws.Range("J20").Value2 = "Gr. 304"
$ws.Range("K20").Value2 = "N/A"
$ws.Range("L20").Value2 = "100 °C"
$ws.Range("M20").Value2 = "1.1 Bar G"
$ws.Range("N20").Value2 = "100 °C"
$ws.Range("O20").Value2 = "1.0 Bar G"
$ws.Range("E21").Value2 = "Earthing Lug"
$ws.Range("G21").Value2 = "DMSO"
$ws.Range("H21").Value2 = "Stainless Steel"
$ws.Range("I21").Value2 = "SA-240"
$ws.Range("J21").Value2 = "Gr. 304"
$ws.Range("K21").Value2 = "N/A"
$ws.Range("L21").Value2 = "100 °C"
$ws.Range("M21").Value2 = "1.1 Bar G"
$ws.Range("N21").Value2 = "100 °C"
$ws.Range("O21").Value2 = "1.0 Bar G"